$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The top three backup codes (rows 2-4) have been used/consumed.
# Shift the remaining unused codes (previously rows 8-10) up into rows 2-4,
# then clear out the now-vacated rows 8-10.
$ws.Range("A2").Value = $ws.Range("A8").Value()
$ws.Range("A3").Value = $ws.Range("A9").Value()
$ws.Range("A4").Value = $ws.Range("A10").Value()

$ws.Range("A8").ClearContents()
$ws.Range("A9").ClearContents()
$ws.Range("A10").ClearContents()

# Move the active selection down to the next row after the refilled codes.
$ws.Range("A5").Select()
